$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.03
$ws.Range("K2").Value = 15

# Row 3
$ws.Range("G3").Value = 4.2
$ws.Range("N3").Value = 1.48
$ws.Range("O3").Value = 2.6

# Row 4
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 2.55
$ws.Range("X4").Value = 26
$ws.Range("AH4").Value = 23
$ws.Range("AI4").Value = 21

# Row 5
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.67

# Row 6
$ws.Range("J6").Value = 1.07
$ws.Range("K6").Value = 9
$ws.Range("L6").Value = 1.36
$ws.Range("M6").Value = 3.2
$ws.Range("N6").Value = 2.1
$ws.Range("O6").Value = 1.73
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.63
$ws.Range("R6").Value = 2.1
$ws.Range("S6").Value = 1.67
$ws.Range("Y6").Value = 34

# Row 8
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("L8").Value = 1.3
$ws.Range("M8").Value = 3.4
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 1.85

# Row 9
$ws.Range("J9").Value = 1.05
$ws.Range("L9").Value = 1.29

# Row 10
$ws.Range("J10").Value = 1.04
$ws.Range("L10").Value = 1.22

# Row 14
$ws.Range("R14").Value = 2.3
$ws.Range("S14").Value = 1.55

# Row 16
$ws.Range("G16").Value = 1.62
$ws.Range("H16").Value = 3.7
$ws.Range("I16").Value = 4.55
$ws.Range("L16").Value = 1.3
$ws.Range("M16").Value = 3.25
$ws.Range("N16").Value = 1.83
$ws.Range("O16").Value = 1.78
$ws.Range("S16").Value = 1.82
$ws.Range("T16").Value = 5.6
$ws.Range("U16").Value = 6.2
$ws.Range("V16").Value = 7
$ws.Range("W16").Value = 9.75
$ws.Range("X16").Value = 11
$ws.Range("Z16").Value = 10
$ws.Range("AA16").Value = 6.3
$ws.Range("AB16").Value = 14
$ws.Range("AC16").Value = 65
$ws.Range("AD16").Value = 450
$ws.Range("AE16").Value = 10.25
$ws.Range("AF16").Value = 21
$ws.Range("AG16").Value = 12.5
$ws.Range("AH16").Value = 60
$ws.Range("AI16").Value = 37
$ws.Range("AJ16").Value = 40

# Row 17
$ws.Range("G17").Value = 2.47
$ws.Range("H17").Value = 3.75
$ws.Range("I17").Value = 2.45
$ws.Range("T17").Value = 11.25
$ws.Range("U17").Value = 14.5
$ws.Range("V17").Value = 9.75
$ws.Range("W17").Value = 27
$ws.Range("X17").Value = 18
$ws.Range("AA17").Value = 7.6
$ws.Range("AE17").Value = 11.25
$ws.Range("AF17").Value = 14
$ws.Range("AG17").Value = 9.5
$ws.Range("AH17").Value = 26
$ws.Range("AI17").Value = 17.5
$ws.Range("AJ17").Value = 23

# Row 19
$ws.Range("G19").Value = 2.07
$ws.Range("H19").Value = 3.45
$ws.Range("I19").Value = 3.2
$ws.Range("L19").Value = 1.22
$ws.Range("M19").Value = 3.5
$ws.Range("O19").Value = 1.98
$ws.Range("R19").Value = 1.55
$ws.Range("S19").Value = 2.15
$ws.Range("T19").Value = 9.25
$ws.Range("U19").Value = 11.25
$ws.Range("W19").Value = 20
$ws.Range("X19").Value = 15
$ws.Range("Y19").Value = 22
$ws.Range("Z19").Value = 12.5
$ws.Range("AA19").Value = 6.9
$ws.Range("AC19").Value = 45
$ws.Range("AE19").Value = 11.75
$ws.Range("AF19").Value = 18.5
$ws.Range("AG19").Value = 11
$ws.Range("AH19").Value = 40
$ws.Range("AI19").Value = 25
$ws.Range("AJ19").Value = 29

# Row 20
$ws.Range("G20").Value = 1.21
$ws.Range("H20").Value = 5.5
$ws.Range("I20").Value = 10
$ws.Range("T20").Value = 7.3
$ws.Range("V20").Value = 7.9
$ws.Range("W20").Value = 6.3
$ws.Range("Z20").Value = 16.5
$ws.Range("AA20").Value = 10
$ws.Range("AE20").Value = 24
$ws.Range("AF20").Value = 60
$ws.Range("AG20").Value = 26
$ws.Range("AH20").Value = 250

# Row 21
$ws.Range("G21").Value = 2.42
$ws.Range("T21").Value = 7.8
$ws.Range("V21").Value = 9.5
$ws.Range("AE21").Value = 8.5
$ws.Range("AF21").Value = 14
$ws.Range("AJ21").Value = 32

# Row 24
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 2.1
$ws.Range("K24").Value = 15
$ws.Range("N24").Value = 1.57
$ws.Range("O24").Value = 2.35
$ws.Range("U24").Value = 19
$ws.Range("X24").Value = 21
$ws.Range("Y24").Value = 23
$ws.Range("AE24").Value = 12
$ws.Range("AF24").Value = 13
$ws.Range("AG24").Value = 9.5
$ws.Range("AH24").Value = 21
$ws.Range("AJ24").Value = 21

# Row 32
$ws.Range("K32").Value = 13

# Row 33
$ws.Range("G33").Value = 3.3
$ws.Range("I33").Value = 2.15
$ws.Range("J33").Value = 1.1
$ws.Range("K33").Value = 7
$ws.Range("T33").Value = 7.5
$ws.Range("V33").Value = 13
$ws.Range("AF33").Value = 9
$ws.Range("AH33").Value = 19
